# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: actualizar mensaje de conversión del día ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newMessage = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.37 = 12984.17 pesos`n✅ 12984.17 pesos = 3.36 = 963.23 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newMessage

# --- tasas: actualizar tasas calculadas ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 296.9
$ws2.Range("O10").Value = 3855
$ws2.Range("N12").Value = 3868.69
$ws2.Range("O12").Value = 287
